# Weekly update: insert two new Tomate price records (Macroferia Regional
# de Talca) ahead of the existing row 488, shifting the remaining rows down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 488-489; everything from the old row 488 on
# shifts down two rows (old 488 -> 490, ... old 520 -> 522).
$ws.Rows("488:489").Insert()

# New row 488
$ws.Cells.Item(488, 1).Value = 5
$ws.Cells.Item(488, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(488, 3).Value = "Maule"
$ws.Cells.Item(488, 4).Value = 44610
$ws.Cells.Item(488, 5).Value = 7
$ws.Cells.Item(488, 6).Value = 100112020
$ws.Cells.Item(488, 7).Value = "Tomate"
$ws.Cells.Item(488, 8).Value = "Larga vida"
$ws.Cells.Item(488, 9).Value = "Primera"
$ws.Cells.Item(488, 10).Value = 2000
$ws.Cells.Item(488, 11).Value = 7000
$ws.Cells.Item(488, 12).Value = 7000
$ws.Cells.Item(488, 13).Value = 7000
$ws.Cells.Item(488, 14).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(488, 15).Value = "Región del Maule"
$ws.Cells.Item(488, 16).Value = 389
$ws.Cells.Item(488, 17).Value = 18
$ws.Cells.Item(488, 18).Value = "Hortaliza"

# New row 489
$ws.Cells.Item(489, 1).Value = 5
$ws.Cells.Item(489, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(489, 3).Value = "Maule"
$ws.Cells.Item(489, 4).Value = 44610
$ws.Cells.Item(489, 5).Value = 7
$ws.Cells.Item(489, 6).Value = 100112020
$ws.Cells.Item(489, 7).Value = "Tomate"
$ws.Cells.Item(489, 8).Value = "Larga vida"
$ws.Cells.Item(489, 9).Value = "Primera"
$ws.Cells.Item(489, 10).Value = 3000
$ws.Cells.Item(489, 11).Value = 4500
$ws.Cells.Item(489, 12).Value = 4500
$ws.Cells.Item(489, 13).Value = 4500
$ws.Cells.Item(489, 14).Value = "`$/caja 15 kilos"
$ws.Cells.Item(489, 15).Value = "Región del Maule"
$ws.Cells.Item(489, 16).Value = 300
$ws.Cells.Item(489, 17).Value = 15
$ws.Cells.Item(489, 18).Value = "Hortaliza"
